# Edit script: insert a new quarterly data column (Period Ending 2018-09-30)
# before the existing "D" column on the ATEC worksheet, shifting the
# previously existing D:K data right to E:L, matching the author's
# "Doing Updates for Financials" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Insert a new column before D; existing D:K (and their formatting)
#    shift one column to the right, to E:L. The freshly inserted column
#    starts out unformatted ("General").
$ws.Columns("D:D").Insert(1)

# 2. The cells that used to be in D now live in E, carrying the original
#    per-row number formats/fonts/alignment (date header rows vs. numeric
#    data rows). Copy that formatting back onto the new (blank) D column
#    so the new quarter's cells look identical to their neighbours.
$ws.Range("E5:E102").Copy() | Out-Null
$ws.Range("D5:D102").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 3. Populate the new D column with the Q3 2018 (period ending 2018-09-30)
#    figures reported for this filing.
$ws.Range("D7").Value = 43373
$ws.Range("D8").Value = 66400
$ws.Range("D9").Value = 19700
$ws.Range("D10").Value = 46700
$ws.Range("D12").Value = 7000
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 300
$ws.Range("D15").Value = 600
$ws.Range("D17").Value = 81100
$ws.Range("D18").Value = -14700
$ws.Range("D20").Value = -5200
$ws.Range("D21").Value = -14900
$ws.Range("D22").Value = "NA"
$ws.Range("D23").Value = -19900
$ws.Range("D24").Value = -1700
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = -18200
$ws.Range("D27").Value = -31700
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = -100
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 5200
$ws.Range("D33").Value = -31800
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -31800
$ws.Range("D38").Value = 43373
$ws.Range("D41").Value = 35100
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 12200
$ws.Range("D44").Value = 29000
$ws.Range("D45").Value = 1800
$ws.Range("D46").Value = 78100
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 12300
$ws.Range("D49").Value = 40800
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 300
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 131500
$ws.Range("D57").Value = 5600
$ws.Range("D58").Value = 6200
$ws.Range("D59").Value = 21300
$ws.Range("D60").Value = 33100
$ws.Range("D61").Value = 34300
$ws.Range("D62").Value = 16200
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 83600
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 23600
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -477800
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 24200
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43373
$ws.Range("D81").Value = -31800
$ws.Range("D83").Value = 5100
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = -16500
$ws.Range("D91").Value = -3300
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -18200
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 47400
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 12600

Write-Output "Done inserting 2018-09-30 quarter column."
